$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild rows 2-10 of the Inhbb -> Acvr1 ligand-receptor pair table.
# The sending/target cluster set grew from {FAPs, sCs} to {ECs, FAPs, sCs},
# so the 2x3 (6-row) matrix becomes a 3x3 (9-row) matrix, and every
# specificity-derived numeric column was recomputed accordingly.

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Inhbb"
$ws.Cells.Item(2,3).Value = "Acvr1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.4201923333333333
$ws.Cells.Item(2,8).Value = 1.260577
$ws.Cells.Item(2,9).Value = 0.08716480679187069
$ws.Cells.Item(2,10).Value = 0.08716480679187069
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 5.057757666666666
$ws.Cells.Item(2,14).Value = 15.173273
$ws.Cells.Item(2,15).Value = 0.173378811020062
$ws.Cells.Item(2,16).Value = 0.173378811020062
$ws.Cells.Item(2,17).Value = 2.125230995391222
$ws.Cells.Item(2,18).Value = 19.127078958521
$ws.Cells.Item(2,19).Value = 0.01511253056436796
$ws.Cells.Item(2,20).Value = 0.01511253056436796
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Inhbb"
$ws.Cells.Item(3,3).Value = "Acvr1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.4201923333333333
$ws.Cells.Item(3,8).Value = 1.260577
$ws.Cells.Item(3,9).Value = 0.08716480679187069
$ws.Cells.Item(3,10).Value = 0.08716480679187069
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 14.247411
$ws.Cells.Item(3,14).Value = 42.742233
$ws.Cells.Item(3,15).Value = 0.4883980890531961
$ws.Cells.Item(3,16).Value = 0.4883980890531961
$ws.Cells.Item(3,17).Value = 5.986652872049
$ws.Cells.Item(3,18).Value = 53.879875848441
$ws.Cells.Item(3,19).Value = 0.0425711250698407
$ws.Cells.Item(3,20).Value = 0.0425711250698407
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Inhbb"
$ws.Cells.Item(4,3).Value = "Acvr1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.4201923333333333
$ws.Cells.Item(4,8).Value = 1.260577
$ws.Cells.Item(4,9).Value = 0.08716480679187069
$ws.Cells.Item(4,10).Value = 0.08716480679187069
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 9.866548666666667
$ws.Cells.Item(4,14).Value = 29.599646
$ws.Cells.Item(4,15).Value = 0.3382230999267418
$ws.Cells.Item(4,16).Value = 0.3382230999267418
$ws.Cells.Item(4,17).Value = 4.145848106193555
$ws.Cells.Item(4,18).Value = 37.312632955742
$ws.Cells.Item(4,19).Value = 0.02948115115766203
$ws.Cells.Item(4,20).Value = 0.02948115115766203
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Inhbb"
$ws.Cells.Item(5,3).Value = "Acvr1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 2.446732
$ws.Cells.Item(5,8).Value = 7.340196000000001
$ws.Cells.Item(5,9).Value = 0.5075507217365239
$ws.Cells.Item(5,10).Value = 0.5075507217365239
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 5.057757666666666
$ws.Cells.Item(5,14).Value = 15.173273
$ws.Cells.Item(5,15).Value = 0.173378811020062
$ws.Cells.Item(5,16).Value = 0.173378811020062
$ws.Cells.Item(5,17).Value = 12.37497753127867
$ws.Cells.Item(5,18).Value = 111.374797781508
$ws.Cells.Item(5,19).Value = 0.08799854066705284
$ws.Cells.Item(5,20).Value = 0.08799854066705286
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Inhbb"
$ws.Cells.Item(6,3).Value = "Acvr1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.446732
$ws.Cells.Item(6,8).Value = 7.340196000000001
$ws.Cells.Item(6,9).Value = 0.5075507217365239
$ws.Cells.Item(6,10).Value = 0.5075507217365239
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 14.247411
$ws.Cells.Item(6,14).Value = 42.742233
$ws.Cells.Item(6,15).Value = 0.4883980890531961
$ws.Cells.Item(6,16).Value = 0.4883980890531961
$ws.Cells.Item(6,17).Value = 34.85959641085201
$ws.Cells.Item(6,18).Value = 313.736367697668
$ws.Cells.Item(6,19).Value = 0.2478868025936888
$ws.Cells.Item(6,20).Value = 0.2478868025936888
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Inhbb"
$ws.Cells.Item(7,3).Value = "Acvr1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.446732
$ws.Cells.Item(7,8).Value = 7.340196000000001
$ws.Cells.Item(7,9).Value = 0.5075507217365239
$ws.Cells.Item(7,10).Value = 0.5075507217365239
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 9.866548666666667
$ws.Cells.Item(7,14).Value = 29.599646
$ws.Cells.Item(7,15).Value = 0.3382230999267418
$ws.Cells.Item(7,16).Value = 0.3382230999267418
$ws.Cells.Item(7,17).Value = 24.14080035229067
$ws.Cells.Item(7,18).Value = 217.267203170616
$ws.Cells.Item(7,19).Value = 0.1716653784757823
$ws.Cells.Item(7,20).Value = 0.1716653784757823
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Inhbb"
$ws.Cells.Item(8,3).Value = "Acvr1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.953740666666667
$ws.Cells.Item(8,8).Value = 5.861222
$ws.Cells.Item(8,9).Value = 0.4052844714716054
$ws.Cells.Item(8,10).Value = 0.4052844714716054
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 5.057757666666666
$ws.Cells.Item(8,14).Value = 15.173273
$ws.Cells.Item(8,15).Value = 0.173378811020062
$ws.Cells.Item(8,16).Value = 0.173378811020062
$ws.Cells.Item(8,17).Value = 9.881546835511777
$ws.Cells.Item(8,18).Value = 88.933921519606
$ws.Cells.Item(8,19).Value = 0.07026773978864118
$ws.Cells.Item(8,20).Value = 0.07026773978864118
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Inhbb"
$ws.Cells.Item(9,3).Value = "Acvr1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.953740666666667
$ws.Cells.Item(9,8).Value = 5.861222
$ws.Cells.Item(9,9).Value = 0.4052844714716054
$ws.Cells.Item(9,10).Value = 0.4052844714716054
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 14.247411
$ws.Cells.Item(9,14).Value = 42.742233
$ws.Cells.Item(9,15).Value = 0.4883980890531961
$ws.Cells.Item(9,16).Value = 0.4883980890531961
$ws.Cells.Item(9,17).Value = 27.835746265414
$ws.Cells.Item(9,18).Value = 250.521716388726
$ws.Cells.Item(9,19).Value = 0.1979401613896667
$ws.Cells.Item(9,20).Value = 0.1979401613896667
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Inhbb"
$ws.Cells.Item(10,3).Value = "Acvr1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.953740666666667
$ws.Cells.Item(10,8).Value = 5.861222
$ws.Cells.Item(10,9).Value = 0.4052844714716054
$ws.Cells.Item(10,10).Value = 0.4052844714716054
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 9.866548666666667
$ws.Cells.Item(10,14).Value = 29.599646
$ws.Cells.Item(10,15).Value = 0.3382230999267418
$ws.Cells.Item(10,16).Value = 0.3382230999267418
$ws.Cells.Item(10,17).Value = 19.27667736971244
$ws.Cells.Item(10,18).Value = 173.490096327412
$ws.Cells.Item(10,19).Value = 0.1370765702932975
$ws.Cells.Item(10,20).Value = 0.1370765702932975
